$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Switch workbook calculation mode to manual (calcPr calcMode="manual")
$excel.Calculation = -4135

# Column D (idtipodni) for data rows 2-74 now holds the text "DNI"
# (shared string) instead of the numeric value 1.
$ws.Range("D2:D74").Value = "DNI"

# Update the sheet's recorded selection/active cell to F54.
$ws.Range("F54").Select()
